$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 5 styling: it becomes a "continuation" row like row 3
#     (A5 gains an empty styled cell, B5:E5 switch to the thin-border style)
$ws.Range("A3:E3").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)  # xlPasteFormats

# --- Add new row 6, formatted like row 4 (the "header" row of a translation pair)
$ws.Range("A4:E4").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(6).RowHeight = 43.2

# Populate the new row's values. Order matches how the shared strings were
# appended in the source workbook (English line, file/line id, translation,
# converted/encoded translation).
$ws.Range("C6").Value = " We made a wonderful discovery!"
$ws.Range("A6").Value = "SCRIPT/P01P04A/us3112.ssb"
$ws.Range("D6").Value = " Мы совершили чудесное открытие!"
$ws.Range("E6").Value = " Íú òïâåñšéìé œôäåòîïå ïóëñúóéå!"
$ws.Range("B6").Value = 71

# --- Restore the selected cell as recorded in the saved workbook
$ws.Range("D2").Select() | Out-Null
